$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resultado (column D) values computed from A (Número A) / B (Operação) / C (Número B)
# using A op C, with "ERRO" for an unrecognised operation or division by zero.
$results = @{
    2  = 31
    3  = -27
    4  = 84
    5  = "ERRO"
    6  = 31
    7  = -19
    8  = 168
    9  = 0.34782608695652201
    10 = 31
    11 = -11
    12 = "ERRO"
    13 = 0.63157894736842102
    14 = 31
    15 = -3
    16 = 240
    17 = "ERRO"
    18 = 31
    19 = 5
    20 = 228
    21 = 1.8181818181818199
    22 = "ERRO"
    23 = 13
    24 = 184
    25 = 3.4285714285714302
    26 = "ERRO"
    27 = 21
    28 = 108
    29 = 9.3333333333333304
    30 = 31
    31 = "ERRO"
}

foreach ($row in 2..31) {
    $ws.Range("D$row").Value = $results[$row]
}
